$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain decimal numbers need to be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# values (losing the "NN.NN" look as text and introducing floating point noise).
$textCells = @("D5","D6","D8","D15","D16","D19","D20","D22","D25","D27","D29","D31","D32","D34","D35","D36","D37","D38","D39","D40","D41","D42","D44","D47","D51")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "60.563.32"
$ws.Range("E2").Value = "  +0.45%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.597.26"
$ws.Range("E3").Value = "  +0.44%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "517.45"
$ws.Range("E5").Value = "  +2.36%  "

# Row 6 (Solana)
$ws.Range("D6").Value = "153.62"
$ws.Range("E6").Value = "  +0.76%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 (XRP)
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  +3.00%  "

# Row 9 (Toncoin)
$ws.Range("E9").Value = "  -0.17%  "

# Row 10 (Dogecoin)
$ws.Range("E10").Value = "  +2.25%  "

# Row 11 (Cardano)
$ws.Range("E11").Value = "  +0.95%  "

# Row 12 (TRON)
$ws.Range("E12").Value = "  +1.80%  "

# Row 13 (WrappedliquidstakedEther2.0)
$ws.Range("D13").Value = "3.051.33"
$ws.Range("E13").Value = "  +0.39%  "

# Row 14 (WrappedBTC)
$ws.Range("D14").Value = "60.519.09"
$ws.Range("E14").Value = "  +0.39%  "

# Row 15 (Avalanche)
$ws.Range("D15").Value = "21.72"
$ws.Range("E15").Value = "  +0.81%  "

# Row 16 (ShibaInu)
$ws.Range("D16").Value = "0.0000141"

# Row 17 (WrappedEther)
$ws.Range("D17").Value = "2.597.04"
$ws.Range("E17").Value = "  +0.27%  "

# Row 18 (Polkadot)
$ws.Range("E18").Value = "  -1.03%  "

# Row 19 (BitcoinCash)
$ws.Range("D19").Value = "351.45"
$ws.Range("E19").Value = "  +1.64%  "

# Row 20 (Chainlink)
$ws.Range("D20").Value = "10.57"
$ws.Range("E20").Value = "  +2.70%  "

# Row 21 (Uniswap)
$ws.Range("E21").Value = "  +2.56%  "

# Row 22 (Dai)
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.11%  "

# Row 23 (Litecoin)
$ws.Range("E23").Value = "  +1.28%  "

# Row 24 (Polygon)
$ws.Range("E24").Value = "  +2.49%  "

# Row 25 and 26 swap (Kaspa moves to row 25, WrappedeETH moves to row 26)
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "0.166"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.712.75"
$ws.Range("E26").Value = "  +0.59%  "

# Row 27 (Binance-PegBSC-USD)
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.37%  "

# Row 28 (PEPE)
$ws.Range("D28").Value = "0.0₃0843"
$ws.Range("E28").Value = "  +0.37%  "

# Row 29 (InternetComputer(DFINITY))
$ws.Range("D29").Value = "7.34"
$ws.Range("E29").Value = "  -1.47%  "

# Row 30 (USDe)
$ws.Range("E30").Value = "  -0.02%  "

# Row 31 (Aptos)
$ws.Range("D31").Value = "6.31"
$ws.Range("E31").Value = "  +10.22%  "

# Row 32 (EthereumClassic)
$ws.Range("D32").Value = "19.40"
$ws.Range("E32").Value = "  +0.80%  "

# Row 33 (PancakeSwap)
$ws.Range("E33").Value = "  +2.67%  "

# Row 34 (Monero)
$ws.Range("D34").Value = "150.42"
$ws.Range("E34").Value = "  -2.92%  "

# Row 35 (NEARProtocol)
$ws.Range("D35").Value = "4.14"
$ws.Range("E35").Value = "  +3.77%  "

# Row 36 (ImmutableX)
$ws.Range("D36").Value = "1.19"
$ws.Range("E36").Value = "  +0.81%  "

# Row 37 (SuiNetwork)
$ws.Range("D37").Value = "0.914"
$ws.Range("E37").Value = "  +6.65%  "

# Row 38 (Stacks)
$ws.Range("D38").Value = "1.49"
$ws.Range("E38").Value = "  +2.78%  "

# Row 39 and 40 swap (Filecoin moves to row 39, OKB moves to row 40)
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "3.78"
$ws.Range("E39").Value = "  +0.21%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "36.36"
$ws.Range("E40").Value = "  +1.66%  "

# Row 41 (Fetch.AI)
$ws.Range("D41").Value = "0.839"
$ws.Range("E41").Value = "  -0.60%  "

# Row 42 (Bittensor)
$ws.Range("D42").Value = "286.78"
$ws.Range("E42").Value = "  -3.68%  "

# Row 43 (Stellar)
$ws.Range("E43").Value = "  +2.07%  "

# Row 44 (Mantle)
$ws.Range("D44").Value = "0.623"
$ws.Range("E44").Value = "  +0.87%  "

# Row 45 (Hedera)
$ws.Range("E45").Value = "  +0.21%  "

# Row 46 (FirstDigitalUSD)
$ws.Range("E46").Value = "  +0.02%  "

# Row 47 (EnergySwap)
$ws.Range("D47").Value = "19.55"
$ws.Range("E47").Value = "  -0.68%  "

# Row 49 (RenderToken)
$ws.Range("E49").Value = "  -0.79%  "

# Row 50 (WhiteBITCoin)
$ws.Range("E50").Value = "  -0.05%  "

# Row 51 (InjectiveProtocol)
$ws.Range("D51").Value = "19.13"
$ws.Range("E51").Value = "  +8.68%  "
